# Updates the cryptos price list (columns D "Price" and E "Volume(1h)")
# for the rows whose figures moved, per the data refresh commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => @{ D = new price text (or $null); E = new volume text (or $null) }
$updates = @{
    2 = @{ D = '60.357.19'; E = '  -3.15%  ' }
    3 = @{ D = '2.896.82'; E = '  -3.70%  ' }
    4 = @{ D = $null; E = '  +0.07%  ' }
    5 = @{ D = '524.31'; E = '  -5.57%  ' }
    6 = @{ D = '140.37'; E = '  -8.15%  ' }
    7 = @{ D = $null; E = '  +0.07%  ' }
    8 = @{ D = '0.548'; E = '  -4.08%  ' }
    9 = @{ D = '2.897.17'; E = '  -3.76%  ' }
    10 = @{ D = '0.106'; E = '  -6.56%  ' }
    11 = @{ D = '5.87'; E = '  -6.56%  ' }
    12 = @{ D = '0.354'; E = '  -3.90%  ' }
    13 = @{ D = '3.402.67'; E = '  -3.73%  ' }
    14 = @{ D = $null; E = '  +0.97%  ' }
    15 = @{ D = '60.490.84'; E = '  -3.16%  ' }
    16 = @{ D = '22.46'; E = '  -5.75%  ' }
    17 = @{ D = '2.903.85'; E = '  -3.58%  ' }
    18 = @{ D = '0.0000139'; E = '  -6.61%  ' }
    19 = @{ D = '4.91'; E = '  -4.26%  ' }
    20 = @{ D = '11.46'; E = '  -4.78%  ' }
    21 = @{ D = '358.50'; E = '  -9.16%  ' }
    22 = @{ D = '6.48'; E = '  -3.46%  ' }
    23 = @{ D = $null; E = '  -0.26%  ' }
    24 = @{ D = '63.22'; E = '  -3.17%  ' }
    25 = @{ D = '3.021.60'; E = '  -3.69%  ' }
    26 = @{ D = '0.445'; E = '  -5.35%  ' }
    27 = @{ D = $null; E = '  -3.44%  ' }
    28 = @{ D = $null; E = '  +0.08%  ' }
    29 = @{ D = '7.73'; E = '  -9.28%  ' }
    30 = @{ D = '0.0₃0838'; E = '  -13.92%  ' }
    31 = @{ D = $null; E = '  -0.01%  ' }
    32 = @{ D = $null; E = '  -5.08%  ' }
    33 = @{ D = '19.35'; E = '  -6.05%  ' }
    34 = @{ D = '150.27'; E = '  -5.93%  ' }
    35 = @{ D = '4.27'; E = '  -9.73%  ' }
    36 = @{ D = '5.49'; E = '  -9.35%  ' }
    37 = @{ D = '0.980'; E = '  -10.29%  ' }
    38 = @{ D = '1.18'; E = '  -9.17%  ' }
    39 = @{ D = '37.77'; E = '  +0.31%  ' }
    40 = @{ D = '1.47'; E = '  -6.99%  ' }
    41 = @{ D = '2.321.79'; E = '  -5.66%  ' }
    42 = @{ D = '0.640'; E = '  -3.42%  ' }
    43 = @{ D = '3.61'; E = '  -8.36%  ' }
    44 = @{ D = '20.53'; E = '  -9.25%  ' }
    45 = @{ D = '0.0568'; E = '  -5.19%  ' }
    46 = @{ D = $null; E = '  +0.01%  ' }
    47 = @{ D = '4.97'; E = '  +0.43%  ' }
    48 = @{ D = '0.0232'; E = '  -7.04%  ' }
    49 = @{ D = '10.33'; E = '  -1.59%  ' }
    50 = @{ D = '0.0922'; E = '  -3.67%  ' }
    51 = @{ D = '248.59'; E = '  -6.34%  ' }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($null -ne $vals.D) {
        # Column D holds plain-text price figures (e.g. "358.50", "0.0₃0838").
        # Force a text number format first so Excel does not reinterpret the
        # string as a numeric value (which would drop significant trailing
        # zeros / normalize thousand-separated figures), then restore the
        # cell to the sheet's normal (unstyled) look, matching column C.
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $vals.D
        $cell.Style = $ws.Range("C$row").Style
    }
    if ($null -ne $vals.E) {
        $ws.Range("E$row").Value = $vals.E
    }
}
